$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'305.49"
$ws.Range("E2").Value = "'-3.86%"
$ws.Range("D3").Value = "'37.20"
$ws.Range("E3").Value = "'-6.55%"
$ws.Range("D4").Value = "'5.092"
$ws.Range("E4").Value = "'-1.13%"
$ws.Range("D5").Value = "'0.07710"
$ws.Range("E5").Value = "'-6.07%"
$ws.Range("D6").Value = "'4.362"
$ws.Range("E6").Value = "'0.52%"
$ws.Range("D7").Value = "'8.209"
$ws.Range("E7").Value = "'-1.64%"
$ws.Range("D8").Value = "'1.867"
$ws.Range("E8").Value = "'-9.74%"
$ws.Range("D9").Value = "'3.193"
$ws.Range("D10").Value = "'0.9184"
$ws.Range("E10").Value = "'-2.35%"
$ws.Range("D11").Value = "'0.1202"
$ws.Range("E11").Value = "'-11.22%"
$ws.Range("D12").Value = "'0.1882"
$ws.Range("E12").Value = "'-5.21%"
$ws.Range("D13").Value = "'0.08700"
$ws.Range("E13").Value = "'-4.48%"
$ws.Range("E14").Value = "'-3.68%"
$ws.Range("D15").Value = "'0.09709"
$ws.Range("E15").Value = "'-1.00%"
$ws.Range("D16").Value = "'0.001367"
$ws.Range("E16").Value = "'-3.99%"
$ws.Range("D17").Value = "'0.006109"
$ws.Range("E17").Value = "'-0.48%"
$ws.Range("D18").Value = "'3.564"
$ws.Range("E18").Value = "'-3.46%"
$ws.Range("D19").Value = "'0.3374"
$ws.Range("E19").Value = "'-2.81%"
$ws.Range("E20").Value = "'-2.75%"
$ws.Range("D21").Value = "'5.028"
$ws.Range("E21").Value = "'1.39%"
$ws.Range("E22").Value = "'1.68%"
$ws.Range("D23").Value = "'0.02113"
$ws.Range("E23").Value = "'5,175.00%"
$ws.Range("D24").Value = "'0.04334"
$ws.Range("E24").Value = "'-0.61%"
$ws.Range("D25").Value = "'0.001216"
$ws.Range("E25").Value = "'-1.21%"
$ws.Range("D26").Value = "'0.004464"
$ws.Range("E26").Value = "'-7.16%"
$ws.Range("D27").Value = "'0.0001356"
$ws.Range("D39").Value = "'0.02223"
$ws.Range("E39").Value = "'-0.68%"
$ws.Range("D40").Value = "'0.04911"
$ws.Range("E40").Value = "'-5.49%"
$ws.Range("D41").Value = "'0.007600"
$ws.Range("E41").Value = "'-1.91%"
$ws.Range("D42").Value = "'0.009911"
$ws.Range("E42").Value = "'2.20%"
$ws.Range("D43").Value = "'0.1328"
$ws.Range("E43").Value = "'-5.39%"
$ws.Range("D44").Value = "'0.002002"
$ws.Range("D45").Value = "'0.008848"
$ws.Range("E45").Value = "'-3.48%"
$ws.Range("D46").Value = "'0.00006791"
$ws.Range("E46").Value = "'2.85%"
$ws.Range("D47").Value = "'0.00000000753"
$ws.Range("E47").Value = "'0.30%"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.001305"
$ws.Range("E48").Value = "'-22.90%"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "'0.003009"
$ws.Range("E49").Value = "'2.11%"
$ws.Range("D50").Value = "'0.00002109"
$ws.Range("E50").Value = "'0.30%"
$ws.Range("D51").Value = "'0.0002009"
$ws.Range("E51").Value = "'0.30%"
